# Add a new 'Strategy' column at column D, shifting the existing
# Return / Return_with_prediction / pct_change / mean_pct_change columns
# one to the right (D:G -> E:H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before the current column D ("Return").
$ws.Range("D1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("D1").Value = "Strategy"

# The data rows alternate between two strategies, starting with
# "Min volatility (Markowitz)" on the first data row (row 2).
$strategies = @("Min volatility (Markowitz)", "Sharpe Ratio")

for ($r = 2; $r -le 29; $r++) {
    $ws.Range("D$r").Value = $strategies[($r) % 2]
}
